$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.230.32"
$ws.Range("E2").Value = "  -0.26%  "

$ws.Range("D3").Value = "3.393.20"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.84"
$ws.Range("E5").Value = "  -0.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.10"
$ws.Range("E6").Value = "  +0.91%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "3.390.91"
$ws.Range("E8").Value = "  +0.34%  "

$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("E10").Value = "  +2.37%  "

$ws.Range("E11").Value = "  -1.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.381"
$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("D13").Value = "3.971.56"
$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("E14").Value = "  -0.83%  "

$ws.Range("E15").Value = "  +2.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").Value = "3.391.87"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").Value = "61.201.14"
$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.88"
$ws.Range("E19").Value = "  +0.35%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.90"
$ws.Range("E20").Value = "  -1.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.34"
$ws.Range("E21").Value = "  -0.37%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.52"
$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.552"
$ws.Range("E23").Value = "  -0.82%  "

$ws.Range("D24").Value = "3.516.23"
$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("E25").Value = "  +0.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.02"
$ws.Range("E26").Value = "  -0.76%  "

$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.64"
$ws.Range("E28").Value = "  -5.47%  "

$ws.Range("E29").Value = "  +8.82%  "

$ws.Range("E30").Value = "  -0.06%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.43"
$ws.Range("E31").Value = "  -0.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("E32").Value = "  -2.20%  "

$ws.Range("E33").Value = "  -1.45%  "

$ws.Range("E34").Value = "  -0.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.48"
$ws.Range("E35").Value = "  +0.14%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.57"
$ws.Range("E36").Value = "  +1.85%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.14"
$ws.Range("E37").Value = "  -2.80%  "

$ws.Range("E38").Value = "  +0.34%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.51"
$ws.Range("E39").Value = "  +0.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0769"
$ws.Range("E40").Value = "  -0.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.07"
$ws.Range("E41").Value = "  +4.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.76"
$ws.Range("E42").Value = "  +2.28%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.775"
$ws.Range("E44").Value = "  +0.09%  "

$ws.Range("E45").Value = "  +1.12%  "

$ws.Range("E46").Value = "  -0.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.18"
$ws.Range("E47").Value = "  -2.84%  "

$ws.Range("D48").Value = "2.531.19"
$ws.Range("E48").Value = "  +7.71%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.62"
$ws.Range("E49").Value = "  +3.84%  "

$ws.Range("E50").Value = "  -0.66%  "

$ws.Range("E51").Value = "  +1.79%  "
